$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Fix the date text on the "第十一周周三" block header
#    (2018.11.12 -> 2018.11.14)
# ------------------------------------------------------------------
$ws.Cells.Item(173, 1).Value = "日期：2018.11.14 第十一周周三"

# ------------------------------------------------------------------
# 2) Append a brand-new weekly block (第十一周四 / 2018.11.15) after
#    row 181, replicating the structure/format of the previous block
#    (rows 163:171) via copy/paste so merged cells, borders and
#    styles match, then overwrite the text content for the new week.
# ------------------------------------------------------------------
$ws.Range("A163:D171").Copy($ws.Range("A185"))

# The engine's merge-paste drops the top-left cell's own border for
# multi-row merges; restore it on the 2-row "总结" block.
$ws.Range("A192:D193").Borders.LineStyle = 1

# Header / date row
$ws.Cells.Item(185, 1).Value = "日期：2018.11.15 第十一周四"

# Data rows (组员 | 计划内容 | 完成情况 | 备注)
$ws.Cells.Item(187, 1).Value = "邱志鹏"
$ws.Cells.Item(187, 2).Value = "更新完善与后台的数据对接"
$ws.Cells.Item(187, 3).Value = "未完成"
$ws.Cells.Item(187, 4).Value = ""

$ws.Cells.Item(188, 1).Value = "黄立根"
$ws.Cells.Item(188, 2).Value = "尝试环信在android端的实现"
$ws.Cells.Item(188, 3).Value = "未完成"
$ws.Cells.Item(188, 4).Value = ""

$ws.Cells.Item(189, 1).Value = "黄俊贤"
$ws.Cells.Item(189, 2).Value = "请假"
$ws.Cells.Item(189, 3).Value = "未完成"
$ws.Cells.Item(189, 4).Value = ""

$ws.Cells.Item(190, 1).Value = "李达波"
$ws.Cells.Item(190, 2).Value = "更新完善与前端的数据对接"
$ws.Cells.Item(190, 3).Value = "未完成"
$ws.Cells.Item(190, 4).Value = ""

$ws.Cells.Item(191, 1).Value = "冯德志"
$ws.Cells.Item(191, 2).Value = "整合两份地图代码，继续地图工作"
$ws.Cells.Item(191, 3).Value = "未完成"
$ws.Cells.Item(191, 4).Value = ""

# "总结：" row (row 192) keeps the text copied from row 170 already.
$ws.Cells.Item(192, 1).Value = "总结："
